$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current totals row (row 9) to make room for
# the new item line (#3), pushing the totals row and the footer row down.
$ws.Rows.Item(9).Insert()

# Copy the formatting from the row above (row 8, item #2) onto the newly
# inserted row 9 so it matches the other item rows (borders, fonts, etc).
$ws.Range("A8:Q8").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row heights: new item row matches the old item-row height, the shifted
# totals row keeps the height row 8 used to have, footer keeps its own.
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75

# Fill in the new item row (#3) values.
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "MAALOX 20 ORAL SACHET SUSP."
$ws.Range("H9").Value = "1:9"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "120.00"
$ws.Range("P9").Value = "18.0000"
$ws.Range("Q9").Value = "0:3"

# Update the running total shown on the (now shifted) totals row.
$ws.Range("P10").Value = 155.5

# Update the footer timestamp on the (now shifted) footer row.
$ws.Range("A11").Value = "Sunday, 28 September, 2025 9:46 AM"
